$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("CoCRoI")
$ws2 = $wb.Worksheets.Item("Overall CAGR")

# --- Sheet: CoCRoI ---
$ws1.Range("B2").Value = -8.201936737888872
$ws1.Range("C2").Value = -3597.574501656507
$ws1.Range("B3").Value = -2.480156882104665
$ws1.Range("C3").Value = -1087.858812413159
$ws1.Range("B4").Value = -2.692992376547527
$ws1.Range("C4").Value = -1181.213781163159
$ws1.Range("B5").Value = -2.585761088550001
$ws1.Range("C5").Value = -1134.179457465244
$ws1.Range("B6").Value = -2.468378947107166
$ws1.Range("C6").Value = -1082.692715674881
$ws1.Range("B7").Value = -2.34062167669367
$ws1.Range("C7").Value = -1026.655182939761
$ws1.Range("B8").Value = -2.202260574728433
$ws1.Range("C8").Value = -965.9665445902588
$ws1.Range("B9").Value = -2.0530624480916
$ws1.Range("C9").Value = -900.5245162941778
$ws1.Range("B10").Value = -1.892789549942089
$ws1.Range("C10").Value = -830.2248163433489
$ws1.Range("B11").Value = -1.721199516918767
$ws1.Range("C11").Value = -754.9611381084942
$ws1.Range("B12").Value = -1.53804530681359
$ws1.Range("C12").Value = -674.6251227011106
$ws1.Range("B13").Value = -1.343075136809851
$ws1.Range("C13").Value = -589.1063318832207
$ws1.Range("B14").Value = -1.136032422384842
$ws1.Range("C14").Value = -498.2922212685514
$ws1.Range("B15").Value = -0.9166557169817396
$ws1.Range("C15").Value = -402.0681138611155
$ws1.Range("B16").Value = -0.68467865256201
$ws1.Range("C16").Value = -300.3171739800116
$ws1.Range("B17").Value = -0.4398298811555597
$ws1.Range("C17").Value = -192.9203816218574
$ws1.Range("B18").Value = -0.1818330175332416
$ws1.Range("C18").Value = -79.75650731551809
$ws1.Range("B19").Value = 0.08959341686695485
$ws1.Range("C19").Value = 39.29791247326807
$ws1.Range("B20").Value = 0.3747360486214085
$ws1.Range("C20").Value = 164.3685993265653
$ws1.Range("B21").Value = 0.6738867054631527
$ws1.Range("C21").Value = 295.5835561837753
$ws1.Range("B22").Value = 0.9873424677305062
$ws1.Range("C22").Value = 433.0730899082932
$ws1.Range("B23").Value = 1.315405717958728
$ws1.Range("C23").Value = 576.9698330396471
$ws1.Range("B24").Value = 1.658384188442401
$ws1.Range("C24").Value = 727.4087646555481
$ws1.Range("B25").Value = 2.016591006587315
$ws1.Range("C25").Value = 884.5272302643612
$ws1.Range("B26").Value = 2.390344737860391
$ws1.Range("C26").Value = 1048.464960644014
$ws1.Range("B27").Value = 2.779969426136218
$ws1.Range("C27").Value = 1219.364089538999
$ws1.Range("B28").Value = 3.185794631228124
$ws1.Range("C28").Value = 1397.369170122436
$ws1.Range("B29").Value = 3.608155463380283
$ws1.Range("C29").Value = 1582.627190125177
$ws1.Range("B30").Value = 4.047392614485933
$ws1.Range("C30").Value = 1775.287585528892
$ws1.Range("B31").Value = 4.503852385784318
$ws1.Range("C31").Value = 1975.502252714647
$ws1.Range("B32").Value = 0.3659616902357496
$ws1.Range("C32").Value = 160.5199463796557

# --- Sheet: Overall CAGR ---
$ws2.Range("B2").Value = -3571.993395746492
$ws2.Range("C2").Value = -3571.993395746492
$ws2.Range("G2").Value = 22156.35913539171
$ws2.Range("B3").Value = -4792.655124826318
$ws2.Range("C3").Value = -1220.661729079825
$ws2.Range("G3").Value = 24761.44194880568
$ws2.Range("H3").Value = -43.5475817639084
$ws2.Range("B4").Value = -5970.282728906143
$ws2.Range("C4").Value = -1177.627604079826
$ws2.Range("G4").Value = 27548.15601071166
$ws2.Range("H4").Value = -20.74994426938493
$ws2.Range("B5").Value = -7100.512354860971
$ws2.Range("C5").Value = -1130.229625954827
$ws2.Range("G5").Value = 30528.20782357103
$ws2.Range("H5").Value = -11.37913675955493
$ws2.Range("B6").Value = -8178.883533237675
$ws2.Range("C6").Value = -1078.371178376704
$ws2.Range("G6").Value = 33713.88387027453
$ws2.Range("H6").Value = -6.367018772471189
$ws2.Range("B7").Value = -9200.837262046804
$ws2.Range("C7").Value = -1021.95372880913
$ws2.Range("B8").Value = -10161.71406261716
$ws2.Range("C8").Value = -960.8768005703614
$ws2.Range("G8").Value = 40754.37519128665
$ws2.Range("H8").Value = -1.217470459994541
$ws2.Range("B9").Value = -11056.75200761351
$ws2.Range("C9").Value = -895.037944996343
$ws2.Range("G9").Value = 44637.00344936474
$ws2.Range("H9").Value = 0.2503621393163291
$ws2.Range("B10").Value = -11881.08472135204
$ws2.Range("C10").Value = -824.3327137385313
$ws2.Range("G10").Value = 48780.96808141219
$ws2.Range("H10").Value = 1.337370072019706
$ws2.Range("B11").Value = -12629.73935258507
$ws2.Range("C11").Value = -748.6546312330363
$ws2.Range("G11").Value = 53202.05451659196
$ws2.Range("H11").Value = 2.16802602528845
$ws2.Range("B12").Value = -13297.6345199651
$ws2.Range("C12").Value = -667.8951673800265
$ws2.Range("G12").Value = 57916.88832642182
$ws2.Range("H12").Value = 2.818481109824966
$ws2.Range("B13").Value = -13879.57823043952
$ws2.Range("C13").Value = -581.9437104744202
$ws2.Range("G13").Value = 62942.98988124645
$ws2.Range("H13").Value = 3.337857574437986
$ws2.Range("B14").Value = -14370.26577087118
$ws2.Range("C14").Value = -490.6875404316543
$ws2.Range("G14").Value = 68298.83287844725
$ws2.Range("H14").Value = 3.759210468235108
$ws2.Range("B15").Value = -14764.27757322588
$ws2.Range("C15").Value = -394.0118023547029
$ws2.Range("G15").Value = 74003.90702461201
$ws2.Range("H15").Value = 4.105564560408226
$ws2.Range("B16").Value = -15056.07705371728
$ws2.Range("C16").Value = -291.7994804914033
$ws2.Range("G16").Value = 80078.78517463114
$ws2.Range("H16").Value = 4.393425776506388
$ws2.Range("B17").Value = -15240.008426351
$ws2.Range("C17").Value = -183.9313726337168
$ws2.Range("B18").Value = -15310.29449136486
$ws2.Range("C18").Value = -70.28606501386366
$ws2.Range("B19").Value = -15261.03439912004
$ws2.Range("C19").Value = 49.2600922448255
$ws2.Range("B20").Value = -15086.20139005877
$ws2.Range("C20").Value = 174.8330090612689
$ws2.Range("B21").Value = -14779.64051140931
$ws2.Range("C21").Value = 306.5608786494577
$ws2.Range("B22").Value = -14335.06631138658
$ws2.Range("C22").Value = 444.5742000227328
$ws2.Range("B23").Value = -13746.06051170903
$ws2.Range("C23").Value = 589.0057996775485
$ws2.Range("B24").Value = -13006.06965932821
$ws2.Range("C24").Value = 739.9908523808208
$ws2.Range("B25").Value = -12108.40275834716
$ws2.Range("C25").Value = 897.6669009810466
$ws2.Range("B26").Value = -11046.22888318833
$ws2.Range("C26").Value = 1062.173875158829
$ws2.Range("B27").Value = -9812.574774160226
$ws2.Range("C27").Value = 1233.654109028105
$ws2.Range("B28").Value = -8400.322416665578
$ws2.Range("C28").Value = 1412.252357494648
$ws2.Range("B29").Value = -6802.206605392172
$ws2.Range("C29").Value = 1598.115811273407
$ws2.Range("B30").Value = -5010.812494930962
$ws2.Range("C30").Value = 1791.394110461209
$ws2.Range("B31").Value = -3018.573138375094
$ws2.Range("C31").Value = 1992.239356555868
$ws2.Range("G31").Value = 233879.108547302
$ws2.Range("H31").Value = 5.941337106727618
